{"js": "// Update the tembakau (tobacco) watering/request interval description and\n// the number of harvests, per commit \"update scoring dan toko\".\n//\n// 1. \"...setiap kelipatan 3 dari umur tembakau...\" -> \"...setiap kelipatan 10 hari dari umur tembakau...\"\n// 2. \"...pada umur tanaman kelipatan 7.\"           -> \"...pada umur tanaman kelipatan 15 hari.\"\n// 3. \"Karena tembakau terdapat 7 kali panen.\"      -> \"Karena tembakau terdapat 3 kali panen.\"\n\nconst body = context.document.body;\n\n// 1) \"kelipatan 3 dari umur tembakau\" -> \"kelipatan 10 hari dari umur tembakau\"\nconst r1 = body.search(\"kelipatan 3 dari umur tembakau\", { matchCase: true });\nr1.load(\"items\");\nawait context.sync();\nif (r1.items.length === 0) {\n  throw new Error(\"Target text for edit #1 not found\");\n}\nr1.items[0].insertText(\"kelipatan 10 hari dari umur tembakau\", Word.InsertLocation.replace);\nawait context.sync();\n\n// 2) \"pada umur tanaman kelipatan 7.\" -> \"pada umur tanaman kelipatan 15 hari.\"\nconst r2 = body.search(\"pada umur tanaman kelipatan 7.\", { matchCase: true });\nr2.load(\"items\");\nawait context.sync();\nif (r2.items.length === 0) {\n  throw new Error(\"Target text for edit #2 not found\");\n}\nr2.items[0].insertText(\"pada umur tanaman kelipatan 15 hari.\", Word.InsertLocation.replace);\nawait context.sync();\n\n// 3) \"Karena tembakau terdapat 7 kali panen.\" -> \"Karena tembakau terdapat 3 kali panen.\"\nconst r3 = body.search(\"Karena tembakau terdapat 7 kali panen.\", { matchCase: true });\nr3.load(\"items\");\nawait context.sync();\nif (r3.items.length === 0) {\n  throw new Error(\"Target text for edit #3 not found\");\n}\nr3.items[0].insertText(\"Karena tembakau terdapat 3 kali panen.\", Word.InsertLocation.replace);\nawait context.sync();\n", "ps1": "# Update the tembakau (tobacco) watering/request interval description and\n# the number of harvests, per commit \"update scoring dan toko\".\n#\n# 1. \"...setiap kelipatan 3 dari umur tembakau...\" -> \"...setiap kelipatan 10 hari dari umur tembakau...\"\n# 2. \"...pada umur tanaman kelipatan 7.\"           -> \"...pada umur tanaman kelipatan 15 hari.\"\n# 3. \"Karena tembakau terdapat 7 kali panen.\"      -> \"Karena tembakau terdapat 3 kali panen.\"\n\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @{Find = \"kelipatan 3 dari umur tembakau\"; Replace = \"kelipatan 10 hari dari umur tembakau\"},\n    @{Find = \"pada umur tanaman kelipatan 7.\"; Replace = \"pada umur tanaman kelipatan 15 hari.\"},\n    @{Find = \"Karena tembakau terdapat 7 kali panen.\"; Replace = \"Karena tembakau terdapat 3 kali panen.\"}\n)\n\nforeach ($r in $replacements) {\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $r.Find\n    $find.Replacement.Text = $r.Replace\n    $find.Forward = $true\n    $find.Wrap = 1\n    $find.MatchCase = $true\n    $find.MatchWholeWord = $false\n    $find.Execute($null, $true, $false, $false, $false, $false, $true, 1, $false, $r.Replace, 2) | Out-Null\n}\n"}
